$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first four data rows (original rows 2-5, years 1985-1988),
# shifting the remaining rows up so the data starts at the original row 6
# (year 1989). This also shrinks the used range from A1:E42 to A1:E38.
$ws.Range("A2:E5").EntireRow.Delete()
